$d = $word.ActiveDocument

$d.Content.Find.Execute("91×65=", $true, $false, $false, $false, $false, $true, 1, $false, "63×16=", 2)
$d.Content.Find.Execute("89×86=", $true, $false, $false, $false, $false, $true, 1, $false, "43×30=", 2)
$d.Content.Find.Execute("79×35=", $true, $false, $false, $false, $false, $true, 1, $false, "51×26=", 2)
$d.Content.Find.Execute("54×59=", $true, $false, $false, $false, $false, $true, 1, $false, "87×36=", 2)
$d.Content.Find.Execute("39×54=", $true, $false, $false, $false, $false, $true, 1, $false, "50×33=", 2)

$d.Content.Find.Execute("89×21=", $true, $false, $false, $false, $false, $true, 1, $false, "80×15=", 2)
$d.Content.Find.Execute("18×42=", $true, $false, $false, $false, $false, $true, 1, $false, "46×97=", 2)
$d.Content.Find.Execute("65×90=", $true, $false, $false, $false, $false, $true, 1, $false, "61×96=", 2)
$d.Content.Find.Execute("52×73=", $true, $false, $false, $false, $false, $true, 1, $false, "34×61=", 2)
$d.Content.Find.Execute("52×75=", $true, $false, $false, $false, $false, $true, 1, $false, "79×47=", 2)

$d.Content.Find.Execute("43×11=", $true, $false, $false, $false, $false, $true, 1, $false, "37×40=", 2)
$d.Content.Find.Execute("86×40=", $true, $false, $false, $false, $false, $true, 1, $false, "17×16=", 2)
$d.Content.Find.Execute("32×30=", $true, $false, $false, $false, $false, $true, 1, $false, "86×97=", 2)
$d.Content.Find.Execute("88×62=", $true, $false, $false, $false, $false, $true, 1, $false, "90×55=", 2)
$d.Content.Find.Execute("25×46=", $true, $false, $false, $false, $false, $true, 1, $false, "93×14=", 2)

$d.Content.Find.Execute("92×43=", $true, $false, $false, $false, $false, $true, 1, $false, "33×55=", 2)
$d.Content.Find.Execute("23×96=", $true, $false, $false, $false, $false, $true, 1, $false, "58×11=", 2)
$d.Content.Find.Execute("84×78=", $true, $false, $false, $false, $false, $true, 1, $false, "35×99=", 2)
$d.Content.Find.Execute("45×90=", $true, $false, $false, $false, $false, $true, 1, $false, "96×22=", 2)
$d.Content.Find.Execute("55×66=", $true, $false, $false, $false, $false, $true, 1, $false, "26×21=", 2)

$d.Content.Find.Execute("97×66=", $true, $false, $false, $false, $false, $true, 1, $false, "22×61=", 2)
$d.Content.Find.Execute("17×80=", $true, $false, $false, $false, $false, $true, 1, $false, "23×29=", 2)
$d.Content.Find.Execute("48×43=", $true, $false, $false, $false, $false, $true, 1, $false, "54×31=", 2)
$d.Content.Find.Execute("62×36=", $true, $false, $false, $false, $false, $true, 1, $false, "13×89=", 2)
$d.Content.Find.Execute("81×66=", $true, $false, $false, $false, $false, $true, 1, $false, "90×44=", 2)
